$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2-7 from 45183 to 45184
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
